$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (simulated "td_sim_1") and D (record_atd / average) corrected values
# Row format: RowNumber, NewC, NewD
$updates = @(
    @(2,  334, 313.5),
    @(3,  75,  59),
    @(4,  1755, 1714),
    @(5,  48,  39.5),
    @(6,  271, 247.5),
    @(7,  43,  65),
    @(8,  131, 72),
    @(9,  77,  60.5),
    @(10, 281, 250.5),
    @(11, 130, 72),
    @(12, 73,  53.5),
    @(13, 179, 156),
    @(14, 117, 109),
    @(15, 50,  41),
    @(16, 102, 94),
    @(17, 57,  57),
    @(18, 144, 118),
    @(19, 122, 118.5),
    @(20, 163, 144.5),
    @(21, 65,  36.5),
    @(22, 60,  35),
    @(23, 69,  47),
    @(24, 230, 237.5),
    @(25, 79,  62.5),
    @(26, 54,  60),
    @(27, 535, 506.5)
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 3).Value = $u[1]
    $ws.Cells.Item($row, 4).Value = $u[2]
}

# Row 28 only has column C (average across corrected C values); no D cell present
$ws.Cells.Item(28, 3).Value = 201.6923076923077
